$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.601.66'
$ws.Range("E2").Value = '  +4.23%  '
$ws.Range("D3").Value = '1.744.54'
$ws.Range("E3").Value = '  +4.48%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.13%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4800'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.65%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2696'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06258'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.48%  '
$ws.Range("D10").Value = '1.743.28'
$ws.Range("E10").Value = '  +4.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07109'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.81'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6165'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.503'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.41'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.93%  '
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = '26.601.10'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006901'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.18%  '
$ws.Range("E20").Value = '  +2.45%  '
$ws.Range("D21").Value = '1.968.82'
$ws.Range("E21").Value = '  +4.42%  '
$ws.Range("E22").Value = '  +4.41%  '
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.348'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.44'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.83%  '
$ws.Range("E27").Value = '  +5.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.423'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '107.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.018'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.762'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07896'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04567'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6393'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9964'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9482'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '113.37'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +17.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.454'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.979'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.003'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("E42").Value = '  +2.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.674'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +16.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3906'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1204'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.712'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.20%  '
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.924'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.260'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3459'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.75%  '
